# Scheduled market-data refresh: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) on leve-profit rows across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 132.5
$ws.Range("I12").Value = 142.66667
$ws.Range("J12").Value = 102
$ws.Range("K12").Value = 142.66667
$ws.Range("L12").Value = 102
$ws.Range("M12").Value = 27.33332999999999
$ws.Range("N12").Value = -442
# Row 98
$ws.Range("H98").Value = 10001
$ws.Range("I98").Value = 10001
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 10001
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()
# Row 105
$ws.Range("H105").Value = 46999.668
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 46999.668
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 46999.668
$ws.Range("N105").Value = -53987.668
# Row 109
$ws.Range("H109").Value = 18821
$ws.Range("I109").Value = 16000
$ws.Range("J109").Value = 19761.334
$ws.Range("K109").Value = 16000
$ws.Range("L109").Value = 19761.334
$ws.Range("M109").Value = -14613
$ws.Range("N109").Value = -22535.334
# Row 114
$ws.Range("H114").Value = 61000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 61000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 61000
$ws.Range("N114").Value = -69678
# Row 122
$ws.Range("H122").Value = 10001
$ws.Range("I122").Value = 10001
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 30003
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
# Row 125
$ws.Range("H125").Value = 600
$ws.Range("I125").Value = 600
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 5400
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -2940
# Row 126
$ws.Range("H126").Value = 34000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 34000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 34000
$ws.Range("N126").Value = -43880
# Row 130
$ws.Range("H130").Value = 200046800
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 200046800
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 200046800
$ws.Range("N130").Value = -200056840
# Row 137
$ws.Range("H137").Value = 1378.8864
$ws.Range("I137").Value = 1182.2963
$ws.Range("J137").Value = 1691.1177
$ws.Range("K137").Value = 3546.8889
$ws.Range("L137").Value = 5073.3531
$ws.Range("M137").Value = -996.8888999999999
$ws.Range("N137").Value = -10173.3531
# Row 139
$ws.Range("H139").Value = 37978.184
$ws.Range("I139").Value = 7000
$ws.Range("J139").Value = 41076
$ws.Range("K139").Value = 7000
$ws.Range("L139").Value = 41076
$ws.Range("M139").Value = -1860
$ws.Range("N139").Value = -51356

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1325.381
$ws.Range("I61").Value = 1021.75
$ws.Range("J61").Value = 2028.5264
$ws.Range("K61").Value = 1021.75
$ws.Range("L61").Value = 2028.5264
$ws.Range("M61").Value = -809.75
$ws.Range("N61").Value = -2452.5264
# Row 74
$ws.Range("H74").Value = 1280.1482
$ws.Range("I74").Value = 1197.0698
$ws.Range("J74").Value = 1604.909
$ws.Range("K74").Value = 1197.0698
$ws.Range("L74").Value = 1604.909
$ws.Range("M74").Value = -323.0698
$ws.Range("N74").Value = -3352.909
# Row 77
$ws.Range("H77").Value = 1280.1482
$ws.Range("I77").Value = 1197.0698
$ws.Range("J77").Value = 1604.909
$ws.Range("K77").Value = 5985.349
$ws.Range("L77").Value = 8024.545
$ws.Range("M77").Value = -1617.349
$ws.Range("N77").Value = -16760.545
# Row 122
$ws.Range("H122").Value = 1724.92
$ws.Range("I122").Value = 1622.3334
$ws.Range("J122").Value = 1878.8
$ws.Range("K122").Value = 4867.0002
$ws.Range("L122").Value = 5636.4
$ws.Range("M122").Value = -2417.0002
$ws.Range("N122").Value = -10536.4
# Row 132
$ws.Range("H132").Value = 3295.9033
$ws.Range("I132").Value = 1284.8636
$ws.Range("J132").Value = 8211.777
$ws.Range("K132").Value = 3854.5908
$ws.Range("L132").Value = 24635.331
$ws.Range("M132").Value = -1324.5908
$ws.Range("N132").Value = -29695.331
# Row 136
$ws.Range("H136").Value = 1325.381
$ws.Range("I136").Value = 1021.75
$ws.Range("J136").Value = 2028.5264
$ws.Range("K136").Value = 3065.25
$ws.Range("L136").Value = 6085.5792
$ws.Range("M136").Value = -515.25
$ws.Range("N136").Value = -11185.5792

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 6300.1177
$ws.Range("I7").Value = 16715
$ws.Range("J7").Value = 619.2727
$ws.Range("K7").Value = 16715
$ws.Range("L7").Value = 619.2727
$ws.Range("M7").Value = -16602
$ws.Range("N7").Value = -845.2727
# Row 31
$ws.Range("H31").Value = 5953981.5
$ws.Range("I31").Value = 1061.0233
$ws.Range("J31").Value = 25644410
$ws.Range("K31").Value = 1061.0233
$ws.Range("L31").Value = 25644410
$ws.Range("M31").Value = -766.0233000000001
$ws.Range("N31").Value = -25645000
# Row 34
$ws.Range("H34").Value = 5953981.5
$ws.Range("I34").Value = 1061.0233
$ws.Range("J34").Value = 25644410
$ws.Range("K34").Value = 1061.0233
$ws.Range("L34").Value = 25644410
$ws.Range("M34").Value = -859.0233000000001
$ws.Range("N34").Value = -25644814
# Row 134
$ws.Range("H134").Value = 1374.9231
$ws.Range("I134").Value = 839.13043
$ws.Range("J134").Value = 2145.125
$ws.Range("K134").Value = 2517.39129
$ws.Range("L134").Value = 6435.375
$ws.Range("M134").Value = 17.60870999999997
$ws.Range("N134").Value = -11505.375

$ws = $wb.Worksheets.Item("CUL")
# Row 82
$ws.Range("H82").Value = 5566.0835
$ws.Range("I82").Value = 1195
$ws.Range("J82").Value = 6440.3
$ws.Range("K82").Value = 3585
$ws.Range("L82").Value = 19320.9
$ws.Range("M82").Value = -3179
$ws.Range("N82").Value = -20132.9
# Row 85
$ws.Range("H85").Value = 5566.0835
$ws.Range("I85").Value = 1195
$ws.Range("J85").Value = 6440.3
$ws.Range("K85").Value = 3585
$ws.Range("L85").Value = 19320.9
$ws.Range("M85").Value = -2181
$ws.Range("N85").Value = -22128.9

$ws = $wb.Worksheets.Item("GSM")
# Row 110
$ws.Range("H110").Value = 58652
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 58652
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 58652
$ws.Range("N110").Value = -66832
# Row 122
$ws.Range("H122").Value = 772470.4
$ws.Range("I122").Value = 1670951.1
$ws.Range("J122").Value = 2344
$ws.Range("K122").Value = 5012853.300000001
$ws.Range("L122").Value = 7032
$ws.Range("M122").Value = -5010403.300000001
$ws.Range("N122").Value = -11932

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
# Row 103
$ws.Range("H103").Value = 34490
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 34490
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 34490
$ws.Range("N103").Value = -36834
# Row 106
$ws.Range("H106").Value = 11608
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 11608
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 11608
$ws.Range("N106").Value = -14132
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 27473.334
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 27473.334
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -28057.334
# Row 55
$ws.Range("H55").Value = 3120
$ws.Range("I55").Value = 2600
$ws.Range("J55").Value = 3900
$ws.Range("K55").Value = 2600
$ws.Range("L55").Value = 3900
$ws.Range("M55").Value = -2323
$ws.Range("N55").Value = -4454
# Row 107
$ws.Range("H107").Value = 653.8182
$ws.Range("I107").Value = 484.7143
$ws.Range("J107").Value = 949.75
$ws.Range("K107").Value = 1454.1429
$ws.Range("L107").Value = 2849.25
$ws.Range("M107").Value = 465.8571000000002
$ws.Range("N107").Value = -6689.25
# Row 109
$ws.Range("H109").Value = 23796.666
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 23796.666
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 23796.666
$ws.Range("N109").Value = -26570.666
